# Update the workbook to match the latest scraped EPEX Spot / CO2 prices.
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "Prix Spot" sheet: column E (18-jun) had no data yet ("-" placeholders).
#    Fill it in with the now-available closing prices.
# ------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Cells.Item(2,5).Value  = 99.3
$wsPrix.Cells.Item(3,5).Value  = 87.44
$wsPrix.Cells.Item(4,5).Value  = 95.03
$wsPrix.Cells.Item(5,5).Value  = 72.94
$wsPrix.Cells.Item(6,5).Value  = 68.3
$wsPrix.Cells.Item(7,5).Value  = 61.53
$wsPrix.Cells.Item(8,5).Value  = 77
$wsPrix.Cells.Item(9,5).Value  = 98.69
$wsPrix.Cells.Item(10,5).Value = 104.56
$wsPrix.Cells.Item(11,5).Value = 60.67
$wsPrix.Cells.Item(12,5).Value = 21.81
$wsPrix.Cells.Item(13,5).Value = 1.17
$wsPrix.Cells.Item(14,5).Value = 0
$wsPrix.Cells.Item(15,5).Value = -0.01
$wsPrix.Cells.Item(16,5).Value = 0
$wsPrix.Cells.Item(17,5).Value = 2.71
$wsPrix.Cells.Item(18,5).Value = 6.27
$wsPrix.Cells.Item(19,5).Value = 52.91
$wsPrix.Cells.Item(20,5).Value = 75.01000000000001
$wsPrix.Cells.Item(21,5).Value = 114.64
$wsPrix.Cells.Item(22,5).Value = 124.9
$wsPrix.Cells.Item(23,5).Value = 111.6
$wsPrix.Cells.Item(24,5).Value = 128.32
$wsPrix.Cells.Item(25,5).Value = 99.01000000000001

# ------------------------------------------------------------------
# 2) "CO2" sheet: a new day (2025-06-16) needs to be inserted as a new
#    row right after the header, pushing the existing 2025-06-17 row down.
# ------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Rows.Item(2).Insert()
# The freshly inserted row inherits the header's bold/bordered style by
# default; strip that back to the plain (unstyled) look used by the rest
# of the data rows.
$wsCo2.Range("A2:B2").ClearFormats()

# Write the new date as literal text (not as an Excel date serial).
$wsCo2.Cells.Item(2,1).NumberFormat = "@"
$wsCo2.Cells.Item(2,1).Value = "2025-06-16"
$wsCo2.Cells.Item(2,2).Value = "-"
# Drop the helper text numbering format so the cell matches the sheet's
# other (unstyled) data cells exactly.
$wsCo2.Range("A2:B2").ClearFormats()
